$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the cryptos data refresh.
# Column D holds price text that can look numeric (e.g. "526.79"),
# so we force a text number format before assigning it and then
# restore the default "Normal" style so no formatting change is
# introduced - this keeps the exact textual representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.753.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.701.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.726.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.51%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.179.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.700.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.845.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0818"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  +8.21%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("E36").Value = "  +8.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.934"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.27%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.872"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.610"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.140.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +1.12%  "
